$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.056089080913442
$ws.Range("D2").Value = 1.061201643151712
$ws.Range("E2").Value = 1.052397825807201
$ws.Range("F2").Value = 1.070209018951906
$ws.Range("I2").Value = 1.048086462909066
$ws.Range("J2").Value = 1.061092539607857
$ws.Range("K2").Value = 1.063926346289904
$ws.Range("L2").Value = 1.05514666977058
$ws.Range("M2").Value = 1.072909467732615
$ws.Range("N2").Value = 1.062599412723115
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.057350517688728
$ws.Range("D3").Value = 1.062247138227951
$ws.Range("E3").Value = 1.053484054537912
$ws.Range("F3").Value = 1.071487314514538
$ws.Range("I3").Value = 1.048486027695289
$ws.Range("J3").Value = 1.062004407836305
$ws.Range("K3").Value = 1.064785833977175
$ws.Range("L3").Value = 1.056045034485927
$ws.Range("M3").Value = 1.07400293790486
$ws.Range("N3").Value = 1.063512575909041
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.058166292447056
$ws.Range("D4").Value = 1.062923206278945
$ws.Range("E4").Value = 1.054186734915481
$ws.Range("F4").Value = 1.072314474007026
$ws.Range("I4").Value = 1.04874310838817
$ws.Range("J4").Value = 1.06259347107561
$ws.Range("K4").Value = 1.065340938875524
$ws.Range("L4").Value = 1.056625560847069
$ws.Range("M4").Value = 1.074709943750582
$ws.Range("N4").Value = 1.064102475685712
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.058509136973292
$ws.Range("D5").Value = 1.063207322246103
$ws.Range("E5").Value = 1.054482099756164
$ws.Range("F5").Value = 1.072662218100135
$ws.Range("I5").Value = 1.048850835360845
$ws.Range("J5").Value = 1.06284088122455
$ws.Range("K5").Value = 1.065574057303775
$ws.Range("L5").Value = 1.056869429984065
$ws.Range("M5").Value = 1.075007041107602
$ws.Range("N5").Value = 1.064350237185441
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.058566695872436
$ws.Range("D6").Value = 1.063255020571736
$ws.Range("E6").Value = 1.054531690397861
$ws.Range("F6").Value = 1.072720606294164
$ws.Range("I6").Value = 1.04886890270842
$ws.Range("J6").Value = 1.062882408920808
$ws.Range("K6").Value = 1.065613184435741
$ws.Range("L6").Value = 1.056910365931999
$ws.Range("M6").Value = 1.07505691761816
$ws.Range("N6").Value = 1.064391823855792
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.058170873970115
$ws.Range("D7").Value = 1.062927003053268
$ws.Range("E7").Value = 1.054190681757886
$ws.Range("F7").Value = 1.07231912055331
$ws.Range("I7").Value = 1.048744549214486
$ws.Range("J7").Value = 1.062596777891715
$ws.Range("K7").Value = 1.065344054786592
$ws.Range("L7").Value = 1.056628820160043
$ws.Range("M7").Value = 1.074713914080113
$ws.Range("N7").Value = 1.064105787197876
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.056515486109912
$ws.Range("D8").Value = 1.061555064603276
$ws.Range("E8").Value = 1.052764960616048
$ws.Range("F8").Value = 1.070641022786264
$ws.Range("I8").Value = 1.048221801427706
$ws.Range("J8").Value = 1.061400912976637
$ws.Range("K8").Value = 1.064217030693717
$ws.Range("L8").Value = 1.05545043787923
$ws.Range("M8").Value = 1.073279124606029
$ws.Range("N8").Value = 1.062908224017451
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.053594841118421
$ws.Range("D9").Value = 1.059134112584839
$ws.Range("E9").Value = 1.050251188136188
$ws.Range("F9").Value = 1.06768402874606
$ws.Range("I9").Value = 1.047289399341699
$ws.Range("J9").Value = 1.059286089521779
$ws.Range("K9").Value = 1.062223025829058
$ws.Range("L9").Value = 1.053367966950805
$ws.Range("M9").Value = 1.070746599009125
$ws.Range("N9").Value = 1.060790397270713
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051645122675175
$ws.Range("D10").Value = 1.057517739912281
$ws.Range("E10").Value = 1.048574245940766
$ws.Range("F10").Value = 1.065712576517219
$ws.Range("I10").Value = 1.046660177468568
$ws.Range("J10").Value = 1.057871015379477
$ws.Range("K10").Value = 1.06088817961473
$ws.Range("L10").Value = 1.051975511698818
$ws.Range("M10").Value = 1.069055269678522
$ws.Range("N10").Value = 1.059373313560817
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050800211760698
$ws.Range("D11").Value = 1.056817238338047
$ws.Range("E11").Value = 1.047847827708227
$ws.Range("F11").Value = 1.064858854422074
$ws.Range("I11").Value = 1.046385897520621
$ws.Range("J11").Value = 1.057257015829935
$ws.Range("K11").Value = 1.060308846519202
$ws.Range("L11").Value = 1.051371558026366
$ws.Range("M11").Value = 1.068322170089228
$ws.Range("N11").Value = 1.058758442061488
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050486270082796
$ws.Range("D12").Value = 1.056556948323287
$ws.Range("E12").Value = 1.047577957696499
$ws.Range("F12").Value = 1.064541730582736
$ws.Range("I12").Value = 1.046283742753225
$ws.Range("J12").Value = 1.057028757012214
$ws.Range("K12").Value = 1.060093453627719
$ws.Range("L12").Value = 1.051147068732183
$ws.Range("M12").Value = 1.068049750086798
$ws.Range("N12").Value = 1.058529859090064
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050553616407333
$ws.Range("D13").Value = 1.056612785637548
$ws.Range("E13").Value = 1.047635847844058
$ws.Range("F13").Value = 1.06460975537744
$ws.Range("I13").Value = 1.046305667759459
$ws.Range("J13").Value = 1.057077728030662
$ws.Range("K13").Value = 1.06013966534412
$ws.Range("L13").Value = 1.051195229450303
$ws.Range("M13").Value = 1.068108190310074
$ws.Range("N13").Value = 1.058578899652976
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050774263388411
$ws.Range("D14").Value = 1.056795724584581
$ws.Range("E14").Value = 1.047825521118645
$ws.Range("F14").Value = 1.0648326411575
$ws.Range("I14").Value = 1.046377458991836
$ws.Range("J14").Value = 1.057238151815761
$ws.Range("K14").Value = 1.060291046224437
$ws.Range("L14").Value = 1.051353004824986
$ws.Range("M14").Value = 1.068299654110863
$ws.Range("N14").Value = 1.058739551258251
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050910197370966
$ws.Range("D15").Value = 1.056908426980175
$ws.Range("E15").Value = 1.047942378930512
$ws.Range("F15").Value = 1.064969966575224
$ws.Range("I15").Value = 1.046421655467557
$ws.Range("J15").Value = 1.057336968688051
$ws.Range("K15").Value = 1.060384290047528
$ws.Range("L15").Value = 1.051450194970984
$ws.Range("M15").Value = 1.068417606071412
$ws.Range("N15").Value = 1.058838508461833
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05170118221691
$ws.Range("D16").Value = 1.05756421703632
$ws.Range("E16").Value = 1.0486224496048
$ws.Range("F16").Value = 1.065769233514567
$ws.Range("I16").Value = 1.046678342025569
$ws.Range("J16").Value = 1.057911737700856
$ws.Range("K16").Value = 1.060926599721157
$ws.Range("L16").Value = 1.052015572630828
$ws.Range("M16").Value = 1.069103907213983
$ws.Range("N16").Value = 1.059414093712563
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.05219716428393
$ws.Range("D17").Value = 1.057975414236879
$ws.Range("E17").Value = 1.049048960468271
$ws.Range("F17").Value = 1.066270571702469
$ws.Range("I17").Value = 1.046838865923653
$ws.Range("J17").Value = 1.05827193534076
$ws.Range("K17").Value = 1.061266417130078
$ws.Range("L17").Value = 1.052369946875161
$ws.Range("M17").Value = 1.069534205099114
$ws.Range("N17").Value = 1.059774802874432
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052486397827992
$ws.Range("D18").Value = 1.05821520071497
$ws.Range("E18").Value = 1.049297709092122
$ws.Range("F18").Value = 1.06656298720191
$ws.Range("I18").Value = 1.046932321008684
$ws.Range("J18").Value = 1.058481910813776
$ws.Range("K18").Value = 1.061464498124921
$ws.Range("L18").Value = 1.052576549935589
$ws.Range("M18").Value = 1.069785119078588
$ws.Range("N18").Value = 1.059985076536704
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.052585008172957
$ws.Range("D19").Value = 1.058296951952796
$ws.Range("E19").Value = 1.049382521269944
$ws.Range("F19").Value = 1.066662692300451
$ws.Range("I19").Value = 1.046964156994919
$ws.Range("J19").Value = 1.058553486447067
$ws.Range("K19").Value = 1.061532016849545
$ws.Range("L19").Value = 1.052646979774384
$ws.Range("M19").Value = 1.069870662211527
$ws.Range("N19").Value = 1.060056753815604
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052143956790268
$ws.Range("D20").Value = 1.057931302673546
$ws.Range("E20").Value = 1.049003202796627
$ws.Range("F20").Value = 1.0662167835626
$ws.Range("I20").Value = 1.046821661420013
$ws.Range("J20").Value = 1.058233302150039
$ws.Range("K20").Value = 1.061229971271394
$ws.Range("L20").Value = 1.052331935967141
$ws.Range("M20").Value = 1.069488045623323
$ws.Range("N20").Value = 1.059736114820149
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050709291254951
$ws.Range("D21").Value = 1.056741856196316
$ws.Range("E21").Value = 1.047769668352447
$ws.Range("F21").Value = 1.064767007244809
$ws.Range("I21").Value = 1.046356325863272
$ws.Range("J21").Value = 1.057190916329281
$ws.Range("K21").Value = 1.060246473942894
$ws.Range("L21").Value = 1.051306548177613
$ws.Range("M21").Value = 1.068243275966097
$ws.Range("N21").Value = 1.058692248691961
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049806655806652
$ws.Range("D22").Value = 1.055993467498566
$ws.Range("E22").Value = 1.046993829424716
$ws.Range("F22").Value = 1.063855396130292
$ws.Range("I22").Value = 1.046062159399923
$ws.Range("J22").Value = 1.056534414734176
$ws.Range("K22").Value = 1.059626935848169
$ws.Range("L22").Value = 1.050660954099248
$ws.Range("M22").Value = 1.06745997802377
$ws.Range("N22").Value = 1.05803481478929
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050285218384108
$ws.Range("D23").Value = 1.056390254086029
$ws.Range("E23").Value = 1.047405142348503
$ws.Range("F23").Value = 1.064338666903781
$ws.Range("I23").Value = 1.04621825384785
$ws.Range("J23").Value = 1.0568825448438
$ws.Range("K23").Value = 1.059955476721759
$ws.Range("L23").Value = 1.051003280985457
$ws.Range("M23").Value = 1.067875282681215
$ws.Range("N23").Value = 1.058383439283601
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052167999165274
$ws.Range("D24").Value = 1.05795123496915
$ws.Range("E24").Value = 1.049023878805197
$ws.Range("F24").Value = 1.066241088128047
$ws.Range("I24").Value = 1.046829435937803
$ws.Range("J24").Value = 1.05825075920204
$ws.Range("K24").Value = 1.061246439983946
$ws.Range("L24").Value = 1.052349111761385
$ws.Range("M24").Value = 1.069508903324416
$ws.Range("N24").Value = 1.059753596663166
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.054350347151185
$ws.Range("D25").Value = 1.059760402864021
$ws.Range("E25").Value = 1.050901243599438
$ws.Range("F25").Value = 1.068448493795013
$ws.Range("I25").Value = 1.047531786814117
$ws.Range("J25").Value = 1.059833729001579
$ws.Range("K25").Value = 1.062739487564774
$ws.Range("L25").Value = 1.053907058181065
$ws.Range("M25").Value = 1.07140183337328
$ws.Range("N25").Value = 1.061338814461389
